$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws.Range("Y3").Value = 1.57
$ws.Range("Z3").Value = 2.25
$ws.Range("AD3").Value = 7.5

# Row 4 updates
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 4.2
$ws.Range("I4").Value = 5.25
$ws.Range("K4").Value = 2.3
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("S4").Value = 1.85
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = 2.4
$ws.Range("V4").Value = 1.55
$ws.Range("W4").Value = 3.25
$ws.Range("X4").Value = 1.33
$ws.Range("Y4").Value = 1.36
$ws.Range("Z4").Value = 3
$ws.Range("AA4").Value = 1.83
$ws.Range("AB4").Value = 1.83
$ws.Range("AI4").Value = 12
$ws.Range("AJ4").Value = 8
$ws.Range("AK4").Value = 17
